# Weekly fruit/vegetable price update: add the newest week's two rows
# (Primera / Segunda) at the top of the data block (row 853) and push
# everything that was there before down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 853-854; existing rows 853.. shift down to 855..
$ws.Rows("853:854").Insert()

# New row 853 - "Primera" quality for the new date (serial 45154)
$ws.Range("A853").Value = 8
$ws.Range("B853").Value = "Terminal La Palmera de La Serena"
$ws.Range("C853").Value = "Coquimbo"
$ws.Range("D853").Value = 45154
$ws.Range("E853").Value = 4
$ws.Range("F853").Value = 100112043
$ws.Range("G853").Value = "Pepino ensalada"
$ws.Range("H853").Value = "Sin especificar"
$ws.Range("I853").Value = "Primera"
$ws.Range("J853").Value = 560
$ws.Range("K853").Value = 10000
$ws.Range("L853").Value = 11000
$ws.Range("M853").Value = 10500
$ws.Range("N853").Value = "`$/caja 60 unidades"
$ws.Range("O853").Value = "Región de Arica y Parinacota"
$ws.Range("P853").Value = 175
$ws.Range("Q853").Value = 60
$ws.Range("R853").Value = "Hortaliza"

# New row 854 - "Segunda" quality for the same new date (serial 45154)
$ws.Range("A854").Value = 8
$ws.Range("B854").Value = "Terminal La Palmera de La Serena"
$ws.Range("C854").Value = "Coquimbo"
$ws.Range("D854").Value = 45154
$ws.Range("E854").Value = 4
$ws.Range("F854").Value = 100112043
$ws.Range("G854").Value = "Pepino ensalada"
$ws.Range("H854").Value = "Sin especificar"
$ws.Range("I854").Value = "Segunda"
$ws.Range("J854").Value = 360
$ws.Range("K854").Value = 6000
$ws.Range("L854").Value = 7000
$ws.Range("M854").Value = 6500
$ws.Range("N854").Value = "`$/caja 80 unidades"
$ws.Range("O854").Value = "Región de Arica y Parinacota"
$ws.Range("P854").Value = 81
$ws.Range("Q854").Value = 80
$ws.Range("R854").Value = "Hortaliza"
